$wb = $excel.ActiveWorkbook

# Fix the casing of the "hunterhall" sheet name to "hunterHall"
$ws = $wb.Worksheets.Item("hunterhall")
$ws.Name = "hunterHall"

# Make "hunterHall" the active/selected sheet (was "stable")
$ws.Activate()
